$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Put the new text into A1 (this mints the shared-strings table entry)
$ws.Range("A1").Value = "Tretyakova inserted important information"

# Leave the selection on A2, like the saved file shows
$ws.Range("A2").Select()
